$d = $word.ActiveDocument

# --- Add the three new character styles (GaNStyle, GaNParagraph, GaNLinks) ---
$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Fix the duplicated trailing date range and apply GaNStyle to the 4 occurrences ---
$oldCampaignText = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 14. – 23. května, 13. – 22. června, 12. – 21. července. Při pozorování použijte hvězdy oblohy, které zobrazujíSouhvězdí Bootes.14. – 23. května, 13. – 22. června, 12. – 21. července"
$newCampaignText = "Informace v této příručce jsou určeny pro pozorovací kampaň probíhající od 14. – 23. května, 13. – 22. června, 12. – 21. července. Při pozorování použijte hvězdy oblohy, které zobrazujíSouhvězdí Bootes."

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Replacement.ClearFormatting()
$rng.Find.Replacement.Style = "GaNStyle"
[void]$rng.Find.Execute($oldCampaignText, $true, $false, $false, $false, $false, $true, 1, $false, $newCampaignText, 2)

# --- Apply GaNLinks style to the "Jeník Hollan..." run, keeping the empty run before it ---
$linksText = "Jeník Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/"
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found = $rng2.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng2.Style = "GaNLinks"
}

Write-Output "done"
